$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 575, pushing the existing 575:586 block down to 578:589.
$ws.Range("A575:A577").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new weekly data (week of 2022-02-03,
# serial date 44595), variety "Tuna", for qualities Extra / Primera / Segunda.

# Row 575 - Tuna / Extra
$ws.Cells.Item(575, 1).Value = 9
$ws.Cells.Item(575, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(575, 3).Value = "Metropolitana"
$ws.Cells.Item(575, 4).Value = 44595
$ws.Cells.Item(575, 5).Value = 13
$ws.Cells.Item(575, 6).Value = 100112027
$ws.Cells.Item(575, 7).Value = "Melón"
$ws.Cells.Item(575, 8).Value = "Tuna"
$ws.Cells.Item(575, 9).Value = "Extra"
$ws.Cells.Item(575, 10).Value = 350
$ws.Cells.Item(575, 11).Value = 700
$ws.Cells.Item(575, 12).Value = 700
$ws.Cells.Item(575, 13).Value = 700
$ws.Cells.Item(575, 14).Value = "$/unidad"
$ws.Cells.Item(575, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(575, 16).Value = 700
$ws.Cells.Item(575, 17).Value = 1
$ws.Cells.Item(575, 18).Value = "Hortaliza"

# Row 576 - Tuna / Primera
$ws.Cells.Item(576, 1).Value = 9
$ws.Cells.Item(576, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(576, 3).Value = "Metropolitana"
$ws.Cells.Item(576, 4).Value = 44595
$ws.Cells.Item(576, 5).Value = 13
$ws.Cells.Item(576, 6).Value = 100112027
$ws.Cells.Item(576, 7).Value = "Melón"
$ws.Cells.Item(576, 8).Value = "Tuna"
$ws.Cells.Item(576, 9).Value = "Primera"
$ws.Cells.Item(576, 10).Value = 600
$ws.Cells.Item(576, 11).Value = 500
$ws.Cells.Item(576, 12).Value = 500
$ws.Cells.Item(576, 13).Value = 500
$ws.Cells.Item(576, 14).Value = "$/unidad"
$ws.Cells.Item(576, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(576, 16).Value = 500
$ws.Cells.Item(576, 17).Value = 1
$ws.Cells.Item(576, 18).Value = "Hortaliza"

# Row 577 - Tuna / Segunda
$ws.Cells.Item(577, 1).Value = 9
$ws.Cells.Item(577, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(577, 3).Value = "Metropolitana"
$ws.Cells.Item(577, 4).Value = 44595
$ws.Cells.Item(577, 5).Value = 13
$ws.Cells.Item(577, 6).Value = 100112027
$ws.Cells.Item(577, 7).Value = "Melón"
$ws.Cells.Item(577, 8).Value = "Tuna"
$ws.Cells.Item(577, 9).Value = "Segunda"
$ws.Cells.Item(577, 10).Value = 400
$ws.Cells.Item(577, 11).Value = 400
$ws.Cells.Item(577, 12).Value = 400
$ws.Cells.Item(577, 13).Value = 400
$ws.Cells.Item(577, 14).Value = "$/unidad"
$ws.Cells.Item(577, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(577, 16).Value = 400
$ws.Cells.Item(577, 17).Value = 1
$ws.Cells.Item(577, 18).Value = "Hortaliza"

Write-Output "applied"
